# Duplicate the last data row (row 3) into two new rows (4 and 5) on both
# worksheets, preserving values and cell formatting (e.g. date number format).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Determine how many columns of data this sheet has (based on row 1 header).
    $lastCol = $ws.Cells.Item(1, $ws.Columns.Count).End(-4159).Column  # xlToLeft

    $srcRow = 3
    $destRows = 4, 5

    foreach ($destRow in $destRows) {
        for ($col = 1; $col -le $lastCol; $col++) {
            $srcCell  = $ws.Cells.Item($srcRow, $col)
            $destCell = $ws.Cells.Item($destRow, $col)

            $destCell.Value2 = $srcCell.Value2

            if ($srcCell.NumberFormat -ne "General") {
                $destCell.NumberFormat = $srcCell.NumberFormat
            }
        }
    }
}
